$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change existing "*" (single asterisk) markers to "**" (double asterisk)
# in column A for rows 7-10 (Qeios related rows being re-marked).
$ws.Range("A7").Value = "**"
$ws.Range("A8").Value = "**"
$ws.Range("A9").Value = "**"
$ws.Range("A10").Value = "**"

# Add new parent-entity markers in column A for rows 12-15.
$ws.Range("A12").Value = "**"
$ws.Range("A13").Value = "*"
$ws.Range("A14").Value = "**"
$ws.Range("A15").Value = "**"

# Update the view state: scroll/freeze so that row 12 is the first
# visible row below the frozen header/columns, and select A15.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("D2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A15").Select()
